$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove old rows 13-25 (they will be rebuilt with corrected content/order)
$ws.Range("A13:C25").EntireRow.Delete()

# Row 13
$ws.Rows.Item(13).RowHeight = 60
$a13 = $ws.Cells.Item(13, 1)
$a13.Value = "Programa resumido:"
$a13.Font.Bold = $true
$a13.VerticalAlignment = -4160
$b13 = $ws.Cells.Item(13, 2)
$b13.Value = "Semestral"
$b13.Font.Bold = $false
$b13.VerticalAlignment = -4160
$b13.WrapText = $true
$c13 = $ws.Cells.Item(13, 3)
$c13.Value = "Semestral"
$c13.Font.Bold = $false
$c13.VerticalAlignment = -4160
$c13.WrapText = $true
$c13.Font.Color = 255

# Row 14
$ws.Rows.Item(14).RowHeight = 60
$a14 = $ws.Cells.Item(14, 1)
$a14.Value = "Short syllabus:"
$a14.Font.Bold = $true
$a14.VerticalAlignment = -4160
$b14 = $ws.Cells.Item(14, 2)
$b14.Value = "Sequences and series, first and second order ordinary differential equations, Solution of equations differential equations by power series, Fourier series and boundary value problems."
$b14.Font.Bold = $false
$b14.VerticalAlignment = -4160
$b14.WrapText = $true
$c14 = $ws.Cells.Item(14, 3)
$c14.Value = "Sequences and series, first and second order ordinary differential equations, Solution of equations differential equations by power series, Fourier series and boundary value problems."
$c14.Font.Bold = $false
$c14.VerticalAlignment = -4160
$c14.WrapText = $true
$c14.Font.Color = 255

# Row 15
$ws.Rows.Item(15).RowHeight = 120
$a15 = $ws.Cells.Item(15, 1)
$a15.Value = "Programa:"
$a15.Font.Bold = $true
$a15.VerticalAlignment = -4160
$b15 = $ws.Cells.Item(15, 2)
$b15.Value = "01/01/2018"
$b15.Font.Bold = $false
$b15.VerticalAlignment = -4160
$b15.WrapText = $true
$c15 = $ws.Cells.Item(15, 3)
$c15.Value = "01/01/2018"
$c15.Font.Bold = $false
$c15.VerticalAlignment = -4160
$c15.WrapText = $true
$c15.Font.Color = 255

# Row 16
$ws.Rows.Item(16).RowHeight = 120
$a16 = $ws.Cells.Item(16, 1)
$a16.Value = "Syllabus:"
$a16.Font.Bold = $true
$a16.VerticalAlignment = -4160
$b16 = $ws.Cells.Item(16, 2)
$b16.Value = "•Sequences and series: Convergence criteria, absolute and conditional convergence, power series, radius convergence, Derivatives and integration term to term.•First and second order ordinary differential equations : Exact and non-exact differential equations, order reduction, Bernulli equation, the method of undetermined coefficients an variations of parameters, solution of differential equations by power series, applications of first and second order differential equations.•Fourier series: Fourier series convergence theorem, Bessel’s Inequality and Parseval’s identity, Partial differential equations and boundary value problems."
$b16.Font.Bold = $false
$b16.VerticalAlignment = -4160
$b16.WrapText = $true
$c16 = $ws.Cells.Item(16, 3)
$c16.Value = "•Sequences and series: Convergence criteria, absolute and conditional convergence, power series, radius convergence, Derivatives and integration term to term.•First and second order ordinary differential equations : Exact and non-exact differential equations, order reduction, Bernulli equation, the method of undetermined coefficients an variations of parameters, solution of differential equations by power series, applications of first and second order differential equations.•Fourier series: Fourier series convergence theorem, Bessel’s Inequality and Parseval’s identity, Partial differential equations and boundary value problems."
$c16.Font.Bold = $false
$c16.VerticalAlignment = -4160
$c16.WrapText = $true
$c16.Font.Color = 255

# Row 17
$a17 = $ws.Cells.Item(17, 1)
$a17.Value = "Avaliação:"
$a17.Font.Bold = $true
$a17.VerticalAlignment = -4160

# Row 18
$ws.Rows.Item(18).RowHeight = 60
$a18 = $ws.Cells.Item(18, 1)
$a18.Value = "Método:"
$a18.Font.Bold = $true
$a18.VerticalAlignment = -4160
$b18 = $ws.Cells.Item(18, 2)
$b18.Value = "6270264 - Juan Fernando Zapata Zapata"
$b18.Font.Bold = $false
$b18.VerticalAlignment = -4160
$b18.WrapText = $true
$c18 = $ws.Cells.Item(18, 3)
$c18.Value = "6270264 - Juan Fernando Zapata Zapata"
$c18.Font.Bold = $false
$c18.VerticalAlignment = -4160
$c18.WrapText = $true
$c18.Font.Color = 255

# Row 19
$ws.Rows.Item(19).RowHeight = 60
$a19 = $ws.Cells.Item(19, 1)
$a19.Value = "Critério:"
$a19.Font.Bold = $true
$a19.VerticalAlignment = -4160
$b19 = $ws.Cells.Item(19, 2)
$b19.Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$b19.Font.Bold = $false
$b19.VerticalAlignment = -4160
$b19.WrapText = $true
$c19 = $ws.Cells.Item(19, 3)
$c19.Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$c19.Font.Bold = $false
$c19.VerticalAlignment = -4160
$c19.WrapText = $true
$c19.Font.Color = 255

# Row 20
$ws.Rows.Item(20).RowHeight = 60
$a20 = $ws.Cells.Item(20, 1)
$a20.Value = "Norma de recuperação:"
$a20.Font.Bold = $true
$a20.VerticalAlignment = -4160
$b20 = $ws.Cells.Item(20, 2)
$b20.Value = "NF≥ 5,0."
$b20.Font.Bold = $false
$b20.VerticalAlignment = -4160
$b20.WrapText = $true
$c20 = $ws.Cells.Item(20, 3)
$c20.Value = "NF≥ 5,0."
$c20.Font.Bold = $false
$c20.VerticalAlignment = -4160
$c20.WrapText = $true
$c20.Font.Color = 255

# Row 21
$ws.Rows.Item(21).RowHeight = 120
$a21 = $ws.Cells.Item(21, 1)
$a21.Value = "Bibliografia:"
$a21.Font.Bold = $true
$a21.VerticalAlignment = -4160
$b21 = $ws.Cells.Item(21, 2)
$b21.Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
$b21.Font.Bold = $false
$b21.VerticalAlignment = -4160
$b21.WrapText = $true
$c21 = $ws.Cells.Item(21, 3)
$c21.Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
$c21.Font.Bold = $false
$c21.VerticalAlignment = -4160
$c21.WrapText = $true
$c21.Font.Color = 255

# Row 22
$a22 = $ws.Cells.Item(22, 1)
$a22.Value = "Requisitos:"
$a22.Font.Bold = $true
$a22.VerticalAlignment = -4160

# Row 23
$ws.Rows.Item(23).RowHeight = 30
$b23 = $ws.Cells.Item(23, 2)
$b23.Value = "LOB1004 -  Cálculo II  (Requisito fraco)`n"
$b23.Font.Bold = $false
$b23.VerticalAlignment = -4160
$b23.WrapText = $true
$c23 = $ws.Cells.Item(23, 3)
$c23.Value = "LOB1004 -  Cálculo II  (Requisito fraco)`n"
$c23.Font.Bold = $false
$c23.VerticalAlignment = -4160
$c23.WrapText = $true
$c23.Font.Color = 255

# Row 24
$ws.Rows.Item(24).RowHeight = 30
$b24 = $ws.Cells.Item(24, 2)
$b24.Value = "LOB1037 -  Àlgebra Linear  (Requisito fraco)`n"
$b24.Font.Bold = $false
$b24.VerticalAlignment = -4160
$b24.WrapText = $true
$c24 = $ws.Cells.Item(24, 3)
$c24.Value = "LOB1037 -  Àlgebra Linear  (Requisito fraco)`n"
$c24.Font.Bold = $false
$c24.VerticalAlignment = -4160
$c24.WrapText = $true
$c24.Font.Color = 255
